$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the contents of C5 (previously held the shared string "empty"),
# leaving the cell's formatting/style intact.
$ws.Range("C5").ClearContents()

# Update the active selection to C5 to match the saved view state.
$ws.Range("C5").Select()
